# Deploying to gh-pages -- add the 2021 column (R) to the indicator table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-RFromQ {
    param(
        [int]$Row,
        $Value
    )
    $src = $ws.Range("Q$Row")
    $dst = $ws.Range("R$Row")

    # Bring over the same number format / font / borders as the Q (2020)
    # column so the new 2021 column matches the rest of the table.
    $src.Copy() | Out-Null
    $dst.PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false

    if ($null -ne $Value) {
        $dst.Value = $Value
    }
}

# Header row: add the 2021 year label next to 2020.
Set-RFromQ 3 2021

# Data rows: the 2021 figures for every indicator/sub-row.
Set-RFromQ 4  0.12641839647678207
Set-RFromQ 5  0.14922981985616976
Set-RFromQ 6  0.10326895933792253
Set-RFromQ 7  3.433011112114915E-2
Set-RFromQ 8  3.6820478077087354E-2
Set-RFromQ 9  3.1930519190242035E-2
Set-RFromQ 10 8.7302929367211068E-2
Set-RFromQ 11 0.10296328329317765
Set-RFromQ 12 7.1859056271889668E-2
Set-RFromQ 13 0.10716050460690947
Set-RFromQ 14 7.9035451351703812E-2
Set-RFromQ 15 0.13553052227085377
Set-RFromQ 16 6.479643687803946E-2
Set-RFromQ 17 7.643825526207898E-2
Set-RFromQ 18 5.3576570965516782E-2
Set-RFromQ 19 5.4163459619715498E-2
Set-RFromQ 20 6.4872252119520635E-2
Set-RFromQ 21 4.3693418784505472E-2
Set-RFromQ 22 5.1373884452794741E-2
Set-RFromQ 23 2.9662368095156877E-2
Set-RFromQ 24 7.2642215296997686E-2
Set-RFromQ 25 0.13772601093442507
Set-RFromQ 26 0.15668565643254884
Set-RFromQ 27 0.11816042869432726
Set-RFromQ 28 0.33417383115107696
Set-RFromQ 29 0.41139191068108794
Set-RFromQ 30 0.24697746624641295
Set-RFromQ 31 0.16773611144997194
Set-RFromQ 32 0.1959922553363346
Set-RFromQ 33 0.13791201213625709
Set-RFromQ 34 $null
Set-RFromQ 35 0
Set-RFromQ 36 0.1
Set-RFromQ 37 0.2

# Move the saved selection to C1 (matches the workbook's last-saved view).
$ws.Range("C1").Select() | Out-Null

Write-Output "Applied 2021 (column R) figures"
